# v1: color image when zoomed in/out
#
# Task Id=38 ("support applying color while zoomed in or out") is complete,
# so it moves from the "Active" sheet to the "Inactive" (Done) sheet, landing
# at the top of that list with today's "Done" date.
#
# Also fixes up the title of the palette-overlap bug (Id=42) to mention that
# it now also covers the status panel.

$wb = $excel.ActiveWorkbook

$activeSheet = $wb.Worksheets.Item("Active")
$inactiveSheet = $wb.Worksheets.Item("Inactive")

# --- 1) Locate the completed task on the Active sheet ------------------
$doneCell = $activeSheet.Columns.Item(1).Find(38)
$doneRow = $doneCell.Row

$taskId = $activeSheet.Cells.Item($doneRow, 1).Value2
$taskTitle = $activeSheet.Cells.Item($doneRow, 2).Value2
$taskCategory = $activeSheet.Cells.Item($doneRow, 4).Value2
$taskCreated = $activeSheet.Cells.Item($doneRow, 5).Value2

# --- 2) Remove it from Active, shifting the remaining rows up ----------
$activeSheet.Rows.Item($doneRow).EntireRow.Delete()

# --- 3) Insert it at the top of the Inactive (Done) list ---------------
$inactiveSheet.Rows.Item(2).EntireRow.Insert()

$inactiveSheet.Range("A2").Value2 = $taskId
$inactiveSheet.Range("B2").Value2 = $taskTitle
$inactiveSheet.Range("C2").Value2 = "Done"
$inactiveSheet.Range("D2").Value2 = $taskCategory
$inactiveSheet.Range("E2").NumberFormat = "@"
$inactiveSheet.Range("E2").Value2 = $taskCreated
$inactiveSheet.Range("F2").NumberFormat = "@"
$inactiveSheet.Range("F2").Value2 = "8/21/2018"

# match the plain "Normal" style used by every other data row
$inactiveSheet.Range("A2:F2").Style = "Normal"

# --- 4) Update the palette-overlap bug title on Active ------------------
$bugCell = $activeSheet.Columns.Item(1).Find(42)
$bugRow = $bugCell.Row
$activeSheet.Cells.Item($bugRow, 2).Value2 = "bug: expanded palette covers part of picturebox and statuspanel"
